# mise à jour section firebase
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("horaire")

# Merge the "React et PWA" / "Netlify" rows into a single row (row 20),
# shift "Formatif React" into row 21, and fix the Examen/Projet Application
# React order for rows 22-23.

$ws.Range("C20").Value = "[React et PWA](pwa.md)<br/>[Netlify](netlify.md)"
$ws.Range("D20").Value = "[Exercice 13 - PWA](exercice13_pwa.md)<br/>[Exercice 14 - Netlify](exercice14_netlify.md)"

$ws.Range("C21").Value = "Formatif React"
$ws.Range("D21").ClearContents()

$ws.Range("C22").Value = "Projet Application React"
$ws.Range("C23").Value = "Examen React"

$ws.Range("C21").Select()
